$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

$ws.Range("D2").Value = "61.865.68"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "2.408.38"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "553.37"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "142.43"
$ws.Range("E6").Value = "  +3.41%  "
$ws.Range("D8").Value = "0.529"
$ws.Range("E8").Value = "  +0.89%  "
$ws.Range("D9").Value = "2.402.44"
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("D10").Value = "0.108"
$ws.Range("E10").Value = "  +1.32%  "
$ws.Range("D11").Value = "0.155"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("D12").Value = "5.37"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "0.350"
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").Value = "25.85"
$ws.Range("E14").Value = "  +3.15%  "
$ws.Range("D15").Value = "0.0000173"
$ws.Range("E15").Value = "  +4.57%  "
$ws.Range("D16").Value = "2.855.38"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").Value = "61.975.69"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").Value = "2.408.66"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("D19").Value = "11.03"
$ws.Range("E19").Value = "  +2.50%  "
$ws.Range("D20").Value = "4.17"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").Value = "321.56"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").Value = "6.68"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "65.17"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("D25").Value = "1.72"
$ws.Range("E25").Value = "  +4.36%  "
$ws.Range("D26").Value = "8.94"
$ws.Range("E26").Value = "  +8.70%  "
$ws.Range("D27").Value = "572.02"
$ws.Range("E27").Value = "  +14.85%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "2.533.62"
$ws.Range("E29").Value = "  +2.49%  "
$ws.Range("D30").Value = "0.0₃0927"
$ws.Range("E30").Value = "  +5.05%  "
$ws.Range("D31").Value = "1.46"
$ws.Range("E31").Value = "  +5.31%  "
$ws.Range("D32").Value = "8.21"
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("D33").Value = "0.147"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("E34").Value = "  +1.93%  "
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  +3.01%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").Value = "5.62"
$ws.Range("E37").Value = "  +4.71%  "
$ws.Range("D38").Value = "4.75"
$ws.Range("E38").Value = "  +1.43%  "
$ws.Range("D39").Value = "0.381"
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("D40").Value = "18.65"
$ws.Range("E40").Value = "  +0.51%  "
$ws.Range("D41").Value = "149.93"
$ws.Range("E41").Value = "  +3.45%  "
$ws.Range("D42").Value = "1.83"
$ws.Range("E42").Value = "  -3.71%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "2.30"
$ws.Range("E44").Value = "  +11.32%  "
$ws.Range("D45").Value = "149.09"
$ws.Range("E45").Value = "  +2.98%  "
$ws.Range("D46").Value = "3.62"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D47").Value = "0.0540"
$ws.Range("E47").Value = "  +3.56%  "
$ws.Range("D48").Value = "20.08"
$ws.Range("E48").Value = "  +4.45%  "
$ws.Range("D49").Value = "0.588"
$ws.Range("E49").Value = "  +2.45%  "
$ws.Range("D50").Value = "0.0920"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("D51").Value = "0.0227"
$ws.Range("E51").Value = "  +2.02%  "
